# The edit removes the "Administrator" bullet point entirely
# (paragraph "Administrator – pristup svim funkcijama softvera, izmena
# informacija o softveru") from the "Korisnici softvera:" list, along
# with its paragraph mark, so the list goes straight from
# "Korisnici softvera:" to "Registrovani korisnik – ...".
#
# (The rest of the underlying XML diff only re-splits/merges runs inside
# the "Funkcije softvera: ..." paragraph without altering any visible
# text, so there is nothing else to change via the document's object
# model.)

$d = $word.ActiveDocument

$needle = "Administrator"

# Walk paragraphs back-to-front so deleting one doesn't disturb the
# indices of paragraphs we still need to examine.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    if ($text -like "*$needle*pristup svim funkcijama softvera*") {
        $para.Range.Delete()
    }
}
